$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 275; existing rows 275-288 shift down to 276-289.
$ws.Rows.Item(275).Insert()

# Fill the new row 275 with the weekly record (same market/category template,
# new date + volume/price figures).
$ws.Range("A275").Value = 4
$ws.Range("B275").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C275").Value = "Los Lagos"
$ws.Range("D275").Value = 44753
$ws.Range("E275").Value = 10
$ws.Range("F275").Value = 100112043
$ws.Range("G275").Value = "Pepino ensalada"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 200
$ws.Range("K275").Value = 23000
$ws.Range("L275").Value = 25000
$ws.Range("M275").Value = 24000
$ws.Range("N275").Value = "$/caja 60 unidades"
$ws.Range("O275").Value = "Región de Arica y Parinacota"
$ws.Range("P275").Value = 400
$ws.Range("Q275").Value = 60
$ws.Range("R275").Value = "Hortaliza"
